# "Finished the multi-document evaluation"
#
# This script:
#  1. Updates the recorded per-week counts on "Sheet1" (both the TA/Mead
#     word-count block in columns N:P and the POI/MP/LP block in columns
#     Z:AE) with the final evaluation numbers, letting the TTEST formulas
#     in R3:T3 and AG3:AI3 recalculate automatically.
#  2. Hides the now-superseded raw columns (B:J) on "Sheet1" since the
#     evaluation is finished and only the summarized columns are needed.
#  3. Makes "Summary" the active/selected sheet (with G8 selected),
#     instead of "Sheet1".

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$data = $wb.Worksheets.Item("Sheet1")

# --- Update evaluation data on Sheet1 -------------------------------------

# Row 3
$data.Range("Z3").Value = 9

# Row 4
$data.Range("N4").Value = 47
$data.Range("O4").Value = 99
$data.Range("P4").Value = 50
$data.Range("AC4").Value = 15
$data.Range("AD4").Value = 33
$data.Range("AE4").Value = 16

# Row 5
$data.Range("N5").Value = 48
$data.Range("O5").Value = 77
$data.Range("P5").Value = 37
$data.Range("AA5").Value = 9
$data.Range("AB5").Value = 12
$data.Range("AC5").Value = 16
$data.Range("AD5").Value = 25
$data.Range("AE5").Value = 12

# Row 6
$data.Range("N6").Value = 55
$data.Range("O6").Value = 53
$data.Range("P6").Value = 59
$data.Range("AC6").Value = 18
$data.Range("AD6").Value = 17
$data.Range("AE6").Value = 19

# Row 7
$data.Range("N7").Value = 38
$data.Range("O7").Value = 71
$data.Range("P7").Value = 25
$data.Range("AA7").Value = 8
$data.Range("AC7").Value = 12
$data.Range("AD7").Value = 23
$data.Range("AE7").Value = 8

# Row 8
$data.Range("N8").Value = 49
$data.Range("O8").Value = 38
$data.Range("P8").Value = 45
$data.Range("AD8").Value = 12

# Row 9
$data.Range("N9").Value = 43
$data.Range("O9").Value = 65
$data.Range("P9").Value = 52
$data.Range("Z9").Value = 8
$data.Range("AC9").Value = 14
$data.Range("AD9").Value = 21
$data.Range("AE9").Value = 17

# Row 10
$data.Range("N10").Value = 30
$data.Range("O10").Value = 67
$data.Range("P10").Value = 34
$data.Range("AA10").Value = 7
$data.Range("AC10").Value = 10
$data.Range("AD10").Value = 22
$data.Range("AE10").Value = 11

# Row 11
$data.Range("N11").Value = 60
$data.Range("O11").Value = 71
$data.Range("P11").Value = 57
$data.Range("Z11").Value = 10
$data.Range("AA11").Value = 7
$data.Range("AC11").Value = 20
$data.Range("AD11").Value = 23
$data.Range("AE11").Value = 19

# Row 12
$data.Range("N12").Value = 37
$data.Range("O12").Value = 38
$data.Range("P12").Value = 53
$data.Range("Z12").Value = 5
$data.Range("AA12").Value = 6
$data.Range("AB12").Value = 11
$data.Range("AC12").Value = 12
$data.Range("AD12").Value = 12
$data.Range("AE12").Value = 17

# Row 13
$data.Range("N13").Value = 57
$data.Range("O13").Value = 43
$data.Range("P13").Value = 77
$data.Range("Z13").Value = 8
$data.Range("AA13").Value = 9
$data.Range("AB13").Value = 7
$data.Range("AC13").Value = 19
$data.Range("AD13").Value = 14
$data.Range("AE13").Value = 25

# Row 14
$data.Range("N14").Value = 48
$data.Range("O14").Value = 46
$data.Range("P14").Value = 47
$data.Range("Z14").Value = 7
$data.Range("AC14").Value = 16
$data.Range("AD14").Value = 15
$data.Range("AE14").Value = 15

# --- Hide the raw-data columns now that the evaluation is finished -------

$data.Range("B1:J1").EntireColumn.Hidden = $true

# --- Switch the active sheet/selection back to the Summary sheet ---------

$summary.Activate()
$summary.Range("G8").Select()
